$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.5 = 30366.94 pesos`n✅ 30366.94 pesos = 7.46 = 942.46 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $text

# --- Update the numeric rate cells on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 133.4
$ws2.Range("O10").Value = 4050.95
$ws2.Range("N12").Value = 4069.5
$ws2.Range("O12").Value = 126.3
